# Remove Emmet Tab short cut
# Locate the keyboard-shortcuts table that contains the "Emmet expand
# abbreviation" row and delete that entire row. Word will re-flow the
# remaining rows upward automatically.

$d = $word.ActiveDocument

for ($i = 1; $i -le $d.Tables.Count; $i++) {
    $t = $d.Tables.Item($i)
    if ($t.Range.Text -like "*Emmet expand abbreviation*") {
        for ($r = 1; $r -le $t.Rows.Count; $r++) {
            $cellText = $t.Cell($r, 2).Range.Text
            if ($cellText -like "*Emmet expand abbreviation*") {
                $t.Rows.Item($r).Delete()
                break
            }
        }
        break
    }
}
